$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range("D2").Value = '31.241.66'
$ws.Range("E2").Value = '  +4.66%  '

# Row 3: Ethereum
$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range("D3").Value = '1.924.77'
$ws.Range("E3").Value = '  +2.56%  '

# Row 4: TetherUSD
$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9959'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.08%  '

# Row 5: BNB
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.86'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.93%  '

# Row 6: USDC
$ws.Range("B6").Value = 'USDC'
$ws.Range("C6").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9966'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.00%  '

# Row 7: XRP
$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4979'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +1.01%  '

# Row 8: Cardano
$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3013'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +3.90%  '

# Row 9: Dogecoin
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06809'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +3.12%  '

# Row 10: WrappedEther
$ws.Range("B10").Value = 'WrappedEther'
$ws.Range("C10").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D10").Value = '1.896.21'
$ws.Range("E10").Value = '  +1.06%  '

# Row 11: Solana
$ws.Range("B11").Value = 'Solana'
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '17.15'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.58%  '

# Row 12: TRON
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07329'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +2.26%  '

# Row 13: Polygon
$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6894'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +3.20%  '

# Row 14: Litecoin
$ws.Range("B14").Value = 'Litecoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '89.32'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +4.43%  '

# Row 15: Polkadot
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.084'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +5.94%  '

# Row 16: WrappedBTC
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '31.078.00'
$ws.Range("E16").Value = '  +4.11%  '

# Row 17: ShibaInu
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008093'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +3.14%  '

# Row 18: Avalanche
$ws.Range("B18").Value = 'Avalanche'
$ws.Range("C18").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.25'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +3.97%  '

# Row 19: Dai
$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9955'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.27%  '

# Row 20: WrappedliquidstakedEther2.0
$ws.Range("B20").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C20").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D20").Value = '2.134.36'
$ws.Range("E20").Value = '  +0.82%  '

# Row 21: BinanceUSD
$ws.Range("B21").Value = 'BinanceUSD'
$ws.Range("C21").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9926'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.37%  '

# Row 22: Uniswap
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.891'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +3.30%  '

# Row 23: BitcoinCash
$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '178.30'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +33.09%  '

# Row 24: Chainlink
$ws.Range("B24").Value = 'Chainlink'
$ws.Range("C24").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.027'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +8.29%  '

# Row 25: Cosmos
$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.377'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +3.13%  '

# Row 26: Monero
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '153.12'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +4.15%  '

# Row 27: EthereumClassic
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.56'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +11.09%  '

# Row 28: LidoDAOToken
$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.963'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.90%  '

# Row 29: Toncoin
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.442'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +4.86%  '

# Row 30: InternetComputer(DFINITY)
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.328'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +3.79%  '

# Row 31: Stellar
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08927'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +3.94%  '

# Row 32: Filecoin
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.108'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +4.90%  '

# Row 33: Hedera
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05334'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +6.38%  '

# Row 34: ImmutableX
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7506'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +6.48%  '

# Row 35: ARBITRUM
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.151'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +3.71%  '

# Row 36: HuobiToken
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.659'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.01%  '

# Row 37: VeChain
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01886'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +15.29%  '

# Row 38: MXToken
$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.758'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +3.52%  '

# Row 39: RenderToken
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.247'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.84%  '

# Row 40: TrustWalletToken
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9473'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +1.97%  '

# Row 41: FraxShare
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.983'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -1.34%  '

# Row 42: TheSandbox
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4387'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +5.18%  '

# Row 43: Quant
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '105.72'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +3.19%  '

# Row 44: Aptos
$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.870'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +4.76%  '

# Row 45: PaxDollar
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.001'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.84%  '

# Row 46: Algorand
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1325'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +5.53%  '

# Row 47: Cronos
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05827'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +2.38%  '

# Row 48: Decentraland
$ws.Range("B48").Value = 'Decentraland'
$ws.Range("C48").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.3928'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +6.17%  '

# Row 49: Elrond
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '33.46'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +2.93%  '

# Row 50: EnergySwap
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.523'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +4.26%  '

# Row 51: NEARProtocol
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.385'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +3.90%  '
